$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - the previous "Valores" sample table is being
# replaced entirely by the real "atendimentos" export.
$ws.Cells.Clear()

# --- Header row -----------------------------------------------------
# Written left to right, except the B2 write is interleaved here (see
# below) so the shared-string table ends up built in the same order as
# the source workbook.
$ws.Range("A1").Value = "Código"
$ws.Range("B1").Value = "Canal"
$ws.Range("C1").Value = "CNPJ"
$ws.Range("D1").Value = "Cliente"

# One of the data rows still carries the legacy "RN-DIREÇÃO" label (from
# before the "UF" column header was renamed), so it is written here,
# between "Cliente" and "UF", to match the shared-string ordering.
$ws.Range("B2").Value = "RN-DIREÇÃO"

$ws.Range("E1").Value = "UF"
$ws.Range("F1").Value = "Produto"
$ws.Range("G1").Value = "Módulo"
$ws.Range("H1").Value = "Data Abertura"
$ws.Range("I1").Value = "Atendente"
$ws.Range("J1").Value = "Problema"
$ws.Range("K1").Value = "Data Solução"
$ws.Range("L1").Value = "Solução"
$ws.Range("M1").Value = "Motivo"
$ws.Range("N1").Value = "Posição"
$ws.Range("O1").Value = "Estágio"
$ws.Range("P1").Value = "Visitas"

# --- Data rows --------------------------------------------------------
$ws.Range("A2").Value = 8349930
$ws.Range("A3").Value = 8354433
$ws.Range("B3").Value = "RN-DIREÇÃO"

# --- Formatting ---------------------------------------------------
# Column A ("Código") is stored as a number but displayed/treated as
# text, with an underlined font - matches the dedicated style the
# attendance codes use. Format A1 directly, then clone that exact
# formatting onto A2:A3 via copy/paste-special so every "Código" cell
# shares the same single cell style.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Font.Underline = $true
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Size column A to fit the (now underlined, text-formatted) code values.
$ws.Columns("A").AutoFit()

# Select the full column A, same as when the code column is picked after
# the import finishes.
$ws.Columns("A").Select() | Out-Null
